# Auto-generated: apply crypto price/volume updates from commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.705.41'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '1.892.78'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = "'244.93"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = "'0.4921"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = "'0.2960"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.98%  '
$ws.Range('D9').Value = "'0.06791"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.63%  '
$ws.Range('D10').Value = '1.890.77'
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('D11').Value = "'17.20"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('D13').Value = "'90.79"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.07%  '
$ws.Range('D14').Value = "'0.6785"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').Value = "'5.041"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').Value = '30.672.71'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').Value = "'0.000007985"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').Value = "'1.000"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = "'13.17"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.80%  '
$ws.Range('D20').Value = '2.132.26'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').Value = "'4.819"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('D23').Value = "'188.27"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +31.74%  '
$ws.Range('D24').Value = "'6.156"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.61%  '
$ws.Range('D25').Value = "'9.373"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('D26').Value = "'155.66"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.08%  '
$ws.Range('E27').Value = '  +12.34%  '
$ws.Range('D28').Value = "'1.901"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = "'1.396"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').Value = "'4.334"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.01%  '
$ws.Range('D31').Value = "'0.09071"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.37%  '
$ws.Range('D32').Value = "'4.012"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = "'0.05208"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.87%  '
$ws.Range('D34').Value = "'0.7500"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.83%  '
$ws.Range('D35').Value = "'1.109"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.33%  '
$ws.Range('D36').Value = "'2.777"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.40%  '
$ws.Range('D37').Value = "'0.01841"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('D38').Value = "'2.689"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').Value = "'2.140"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').Value = "'0.9379"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('D41').Value = "'0.4419"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.37%  '
$ws.Range('D42').Value = "'105.27"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.90%  '
$ws.Range('D43').Value = "'1.001"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').Value = "'5.762"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').Value = "'7.596"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.73%  '
$ws.Range('D46').Value = "'0.1342"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.18%  '
$ws.Range('E47').Value = '  +2.88%  '
$ws.Range('D48').Value = "'8.727"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.08%  '
$ws.Range('D49').Value = "'1.423"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.83%  '
$ws.Range('D50').Value = "'0.3937"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.96%  '
$ws.Range('E51').Value = '  +2.32%  '
